$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: Guest / read stories / I can read stories without signing up / P1 / Sprint 2 / Done ---
# Copy formatting from existing rows/cells that already carry the exact
# style indices required for the new row, then overwrite the values.
$ws.Range("B33:E33").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F34").PasteSpecial(-4122)
$ws.Range("G12").Copy()
$ws.Range("G34").PasteSpecial(-4122)

$ws.Range("B34").Value = "Guest"
$ws.Range("C34").Value = "read stories "
$ws.Range("D34").Value = "I can read stories without signing up "
$ws.Range("E34").Value = "P1"
$ws.Range("F34").Value = 2
$ws.Range("G34").Value = "Done"

# --- Row 35: User / browse in multi-languages / I can understand everything / P1 / to be started ---
$ws.Range("B33:E33").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("G33").Copy()
$ws.Range("G35").PasteSpecial(-4122)

$ws.Range("B35").Value = "User"
$ws.Range("C35").Value = "browse in multi-languages"
$ws.Range("D35").Value = "I can understand everything "
$ws.Range("E35").Value = "P1"
$ws.Range("G35").Value = "to be started"

$ws.Range("G35").Select()
